# Auto-generated edit script applying the Ifrit_Profits.xlsx diff
# (sheet names ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR correspond to the
# workbook's physical sheet1..sheet8.xml in this file)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 95
$ws.Range("H95").Value = 45000  # was 42500
$ws.Range("J95").Value = 45000  # was 42500
$ws.Range("L95").Value = 45000  # was 42500
$ws.Range("N95").Value = -50492  # was -47992

# Row 112
$ws.Range("H112").Value = 62501252  # was 41667924
$ws.Range("I112").Value = 607.5  # was 566
$ws.Range("J112").Value = 83334800  # was 52633020
$ws.Range("K112").Value = 1822.5  # was 1698
$ws.Range("L112").Value = 250004400  # was 157899060
$ws.Range("M112").Value = -714.5  # was -590
$ws.Range("N112").Value = -250006616  # was -157901276

# Row 132
$ws.Range("H132").Value = 234998.89  # was 235080.8
$ws.Range("I132").Value = 252473.8  # was 265751.94
$ws.Range("J132").Value = 2000  # was 1980
$ws.Range("K132").Value = 757421.3999999999  # was 797255.8200000001
$ws.Range("L132").Value = 6000  # was 5940
$ws.Range("M132").Value = -754891.3999999999  # was -794725.8200000001
$ws.Range("N132").Value = -11060  # was -11000

# Row 138
$ws.Range("H138").Value = 2513.347  # was 2478.848
$ws.Range("I138").Value = 2222.6086  # was 2153.875
$ws.Range("J138").Value = 2770.5386  # was 2833.3635
$ws.Range("K138").Value = 6667.825800000001  # was 6461.625
$ws.Range("L138").Value = 8311.6158  # was 8500.0905
$ws.Range("M138").Value = -1527.825800000001  # was -1321.625
$ws.Range("N138").Value = -18591.6158  # was -18780.0905

# Row 139
$ws.Range("H139").Value = 67500  # was 54500
$ws.Range("J139").Value = 67500  # was 54500
$ws.Range("L139").Value = 67500  # was 54500
$ws.Range("N139").Value = -77780  # was -64780

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4374.6  # was 5411.604
$ws.Range("I32").Value = 4641.087  # was 6098.1797
$ws.Range("K32").Value = 4641.087  # was 6098.1797
$ws.Range("M32").Value = -4354.087  # was -5811.1797

# Row 74
$ws.Range("H74").Value = 3896.5952  # was 3433.0444
$ws.Range("I74").Value = 940  # was 805.6429000000001
$ws.Range("J74").Value = 4945.7095  # was 4619.613
$ws.Range("K74").Value = 940  # was 805.6429000000001
$ws.Range("L74").Value = 4945.7095  # was 4619.613
$ws.Range("M74").Value = -66  # was 68.35709999999995
$ws.Range("N74").Value = -6693.7095  # was -6367.613

# Row 77
$ws.Range("H77").Value = 3896.5952  # was 3433.0444
$ws.Range("I77").Value = 940  # was 805.6429000000001
$ws.Range("J77").Value = 4945.7095  # was 4619.613
$ws.Range("K77").Value = 4700  # was 4028.2145
$ws.Range("L77").Value = 24728.5475  # was 23098.065
$ws.Range("M77").Value = -332  # was 339.7855
$ws.Range("N77").Value = -33464.5475  # was -31834.065

# Row 92
$ws.Range("H92").Value = 20550  # was 0
$ws.Range("J92").Value = 20550  # was 0
$ws.Range("L92").Value = 20550  # was 0
$ws.Range("N92").Value = -25542

# Row 106
$ws.Range("H106").Value = 35000  # was 35333.332
$ws.Range("J106").Value = 35000  # was 35333.332
$ws.Range("L106").Value = 35000  # was 35333.332
$ws.Range("N106").Value = -37524  # was -37857.332

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 910.1  # was 1318.4166
$ws.Range("I99").Value = 996.6667  # was 1542
$ws.Range("J99").Value = 873  # was 1158.7142
$ws.Range("K99").Value = 996.6667  # was 1542
$ws.Range("L99").Value = 873  # was 1158.7142
$ws.Range("M99").Value = 501.3333  # was -44
$ws.Range("N99").Value = -3869  # was -4154.7142

# Row 103
$ws.Range("H103").Value = 24494.5  # was 28999.5
$ws.Range("J103").Value = 24494.5  # was 28999.5
$ws.Range("L103").Value = 24494.5  # was 28999.5
$ws.Range("N103").Value = -26838.5  # was -31343.5

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1297.2  # was 1227.5555
$ws.Range("I99").Value = 1311.3334  # was 1178.2858
$ws.Range("J99").Value = 1276  # was 1400
$ws.Range("K99").Value = 1311.3334  # was 1178.2858
$ws.Range("L99").Value = 1276  # was 1400
$ws.Range("M99").Value = 186.6666  # was 319.7141999999999
$ws.Range("N99").Value = -4272  # was -4396

# Row 107
$ws.Range("H107").Value = 404.4  # was 412.88
$ws.Range("I107").Value = 395.9  # was 399.2
$ws.Range("J107").Value = 438.4  # was 467.6
$ws.Range("K107").Value = 395.9  # was 399.2
$ws.Range("L107").Value = 438.4  # was 467.6
$ws.Range("M107").Value = 1524.1  # was 1520.8
$ws.Range("N107").Value = -4278.4  # was -4307.6

# Row 122
$ws.Range("H122").Value = 4167850.2  # was 4311571
$ws.Range("I122").Value = 6579956  # was 6945512
$ws.Range("K122").Value = 19739868  # was 20836536
$ws.Range("M122").Value = -19737418  # was -20834086

# Row 126
$ws.Range("H126").Value = 1297.2  # was 1227.5555
$ws.Range("I126").Value = 1311.3334  # was 1178.2858
$ws.Range("J126").Value = 1276  # was 1400
$ws.Range("K126").Value = 3934.0002  # was 3534.8574
$ws.Range("L126").Value = 3828  # was 4200
$ws.Range("M126").Value = -1464.0002  # was -1064.8574
$ws.Range("N126").Value = -8768  # was -9140

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 669.4  # was 219.6
$ws.Range("I13").Value = 474  # was 219.6
$ws.Range("J13").Value = 799.6667  # was 0
$ws.Range("K13").Value = 1422  # was 658.8
$ws.Range("L13").Value = 2399.0001  # was 0
$ws.Range("M13").Value = -1254  # was -490.8
$ws.Range("N13").Value = -2735.0001

# Row 131
$ws.Range("H131").Value = 2384408.5  # was 3230108
$ws.Range("J131").Value = 2859514  # was 4169632.2
$ws.Range("L131").Value = 8578542  # was 12508896.6
$ws.Range("N131").Value = -8588622  # was -12518976.6

# Row 137
$ws.Range("H137").Value = 37511.2  # was 30522.719
$ws.Range("I137").Value = 2051.818  # was 2056.3635
$ws.Range("J137").Value = 58040.316  # was 41705.93
$ws.Range("K137").Value = 6155.454000000001  # was 6169.0905
$ws.Range("L137").Value = 174120.948  # was 125117.79
$ws.Range("M137").Value = -1055.454000000001  # was -1069.0905
$ws.Range("N137").Value = -184320.948  # was -135317.79

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 2183.446  # was 2130.726
$ws.Range("I5").Value = 0  # was 200
$ws.Range("J5").Value = 2183.446  # was 2157.5417
$ws.Range("K5").Value = 0  # was 200
$ws.Range("L5").Value = 2183.446  # was 2157.5417
$ws.Range("M5").ClearContents()  # cell removed (was -88)
$ws.Range("N5").Value = -2407.446  # was -2381.5417

# Row 22
$ws.Range("H22").Value = 900  # was 10000
$ws.Range("J22").Value = 900  # was 10000
$ws.Range("L22").Value = 900  # was 10000
$ws.Range("N22").Value = -1958  # was -11058

# Row 95
$ws.Range("H95").Value = 18671.75  # was 26258
$ws.Range("J95").Value = 18671.75  # was 26258
$ws.Range("L95").Value = 18671.75  # was 26258
$ws.Range("N95").Value = -24163.75  # was -31750

# Row 102
$ws.Range("H102").Value = 1383.4  # was 1332.6875
$ws.Range("I102").Value = 1002.875  # was 940.9
$ws.Range("J102").Value = 1818.2858  # was 1985.6666
$ws.Range("K102").Value = 1002.875  # was 940.9
$ws.Range("L102").Value = 1818.2858  # was 1985.6666
$ws.Range("M102").Value = 619.125  # was 681.1
$ws.Range("N102").Value = -5062.2858  # was -5229.6666

# Row 113
$ws.Range("H113").Value = 2367.7273  # was 2310.6667
$ws.Range("I113").Value = 2947.8  # was 2290.7144
$ws.Range("J113").Value = 1884.3334  # was 2338.6
$ws.Range("K113").Value = 2947.8  # was 2290.7144
$ws.Range("L113").Value = 1884.3334  # was 2338.6
$ws.Range("M113").Value = -777.8000000000002  # was -120.7143999999998
$ws.Range("N113").Value = -6224.3334  # was -6678.6

# Row 126
$ws.Range("H126").Value = 3720.2  # was 1292.2222
$ws.Range("I126").Value = 3775.25  # was 1328.75
$ws.Range("J126").Value = 3500  # was 1000
$ws.Range("K126").Value = 11325.75  # was 3986.25
$ws.Range("L126").Value = 10500  # was 3000
$ws.Range("M126").Value = -8855.75  # was -1516.25
$ws.Range("N126").Value = -15440  # was -7940

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 5171992  # was 7754000.5
$ws.Range("J2").Value = 6006390.5  # was 10005334
$ws.Range("L2").Value = 6006390.5  # was 10005334
$ws.Range("N2").Value = -6006614.5  # was -10005558

# Row 7
$ws.Range("H7").Value = 1895.3846  # was 1839.2812
$ws.Range("I7").Value = 1804.375  # was 1716.6086
$ws.Range("J7").Value = 2041  # was 2152.7778
$ws.Range("K7").Value = 1804.375  # was 1716.6086
$ws.Range("L7").Value = 2041  # was 2152.7778
$ws.Range("M7").Value = -1692.375  # was -1604.6086
$ws.Range("N7").Value = -2265  # was -2376.7778

# Row 40
$ws.Range("H40").Value = 2337.4375  # was 2555.5557
$ws.Range("I40").Value = 2254.4546  # was 2475
$ws.Range("J40").Value = 2520  # was 3200
$ws.Range("K40").Value = 2254.4546  # was 2475
$ws.Range("L40").Value = 2520  # was 3200
$ws.Range("M40").Value = -2118.4546  # was -2339
$ws.Range("N40").Value = -2792  # was -3472

# Row 45
$ws.Range("H45").Value = 14748.5  # was 14448.5
$ws.Range("J45").Value = 14998  # was 14398
$ws.Range("L45").Value = 14998  # was 14398
$ws.Range("N45").Value = -15812  # was -15212

# Row 48
$ws.Range("H48").Value = 12000  # was 11950
$ws.Range("J48").Value = 14000  # was 13900
$ws.Range("L48").Value = 14000  # was 13900
$ws.Range("N48").Value = -15322  # was -15222

# Row 62
$ws.Range("H62").Value = 16890  # was 22299.666
$ws.Range("J62").Value = 16890  # was 22299.666
$ws.Range("L62").Value = 16890  # was 22299.666
$ws.Range("N62").Value = -18138  # was -23547.666

# Row 65
$ws.Range("H65").Value = 16890  # was 22299.666
$ws.Range("J65").Value = 16890  # was 22299.666
$ws.Range("L65").Value = 50670  # was 66898.99800000001
$ws.Range("N65").Value = -56910  # was -73138.99800000001

# Row 95
$ws.Range("H95").Value = 14528.8  # was 16666.334
$ws.Range("J95").Value = 14528.8  # was 16666.334
$ws.Range("L95").Value = 14528.8  # was 16666.334
$ws.Range("N95").Value = -20020.8  # was -22158.334

# Row 98
$ws.Range("H98").Value = 31674  # was 33354
$ws.Range("J98").Value = 31674  # was 33354
$ws.Range("L98").Value = 31674  # was 33354
$ws.Range("N98").Value = -37664  # was -39344

# Row 126
$ws.Range("H126").Value = 1895.3846  # was 1839.2812
$ws.Range("I126").Value = 1804.375  # was 1716.6086
$ws.Range("J126").Value = 2041  # was 2152.7778
$ws.Range("K126").Value = 5413.125  # was 5149.825800000001
$ws.Range("L126").Value = 6123  # was 6458.3334
$ws.Range("M126").Value = -2943.125  # was -2679.825800000001
$ws.Range("N126").Value = -11063  # was -11398.3334

# Row 132
$ws.Range("H132").Value = 8136.143  # was 189183.67
$ws.Range("I132").Value = 12673.272  # was 280776
$ws.Range("J132").Value = 3145.3  # was 5999
$ws.Range("K132").Value = 38019.81600000001  # was 842328
$ws.Range("L132").Value = 9435.900000000001  # was 17997
$ws.Range("M132").Value = -35489.81600000001  # was -839798
$ws.Range("N132").Value = -14495.9  # was -23057

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 3000  # was 230136.36
$ws.Range("J18").Value = 3000  # was 230136.36
$ws.Range("L18").Value = 3000  # was 230136.36
$ws.Range("N18").Value = -3346  # was -230482.36

# Row 136
$ws.Range("H136").Value = 8040.7856  # was 7330.516
$ws.Range("I136").Value = 11824.556  # was 11218.211
$ws.Range("J136").Value = 1230  # was 1175
$ws.Range("K136").Value = 35473.66800000001  # was 33654.633
$ws.Range("L136").Value = 3690  # was 3525
$ws.Range("M136").Value = -32923.66800000001  # was -31104.633
$ws.Range("N136").Value = -8790  # was -8625
